$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where Species (F) = "Na" -> harvest record should be cleared/zeroed out.
$naRows = @(3,4,6,8,9,10,11,12,15,16,17,18,19,20,24,27,28,29,30,31,32,33,34,35,37,38,39,40,41,42,45,46,47,49,52,53,57,58,59,61,65,67,68,69,70,73,74,75,76,81)

# Subset of those rows that also had "Harvest in Park?" = Yes, which must be
# reset along with Park Name / Permit Authorization Number.
$parkRows = @(8,9,10,11,12,15,16,17,18,19,20,32,33)

foreach ($r in $naRows) {
    $ws.Cells.Item($r, 2).Value = "No"     # B: Did Harvest Occur?
    $ws.Cells.Item($r, 6).Value = ""       # F: Species
    $ws.Cells.Item($r, 10).Value = 0       # J: Unknown Sex Count
}

foreach ($r in $parkRows) {
    $ws.Cells.Item($r, 11).Value = "No"    # K: Harvest in Park?
    $ws.Cells.Item($r, 12).Value = ""      # L: Park Name
    $ws.Cells.Item($r, 13).Value = ""      # M: PERMITAUTHORIZATIONNUMBER
}
